# Jeannine's log.xlsx update
# - Add a new Wednesday (2016-08-24) AV Shutdown / LSB 105 entry to the "Logs" sheet
# - Apply the new row height used on recent entries to the last block of rows
# - Rebuild the "database" lookup sheet: drop the "Crestron Logout" task type and
#   move Carl / ATK to the end of their respective lookup lists

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$db   = $wb.Worksheets.Item("database")

# ---------------------------------------------------------------------------
# 1. "Logs" sheet - refresh the row height used on the newest block of entries
# ---------------------------------------------------------------------------
$logs.Rows("50").RowHeight = 14.45
$logs.Rows("51").RowHeight = 14.45
$logs.Rows("52").RowHeight = 14.45
$logs.Rows("53").RowHeight = 14.45

# ---------------------------------------------------------------------------
# 2. "Logs" sheet - append a new WEDNESDAY group header + log entry
# ---------------------------------------------------------------------------
$logs.Range("A49:F49").Copy($logs.Range("A56:F56"))
$logs.Range("B56").Value = "WEDNESDAY"

$logs.Range("A50:F50").Copy($logs.Range("A57:F57"))
$logs.Range("A57").Value = "AV Shutdown"
$logs.Range("B57").Value = 42606
$logs.Range("C57").Value = "1730"
$logs.Range("D57").Value = "LSB"
$logs.Range("E57").Value = "105"
$logs.Range("F57").Value = "Log off and make sure neck mic goes back to drawer."

# ---------------------------------------------------------------------------
# 3. "database" sheet - Staff_Name (A) / Task_type (B) / Building (C) lists
#    Carl (A4/C4 pair "ATK") move to the bottom of their lists; the
#    "Crestron Logout" task type is removed outright (not re-added).
# ---------------------------------------------------------------------------
$staff = @("Christina","Christine","CindyM","CindyT","Clairissa","Daniel","Dean","Eli","Elizabeth","Eric","Filipe","Ghazal","Hasebullah","Hashir","Jackie","Jasleen","Jhan","Keagan","Konrad","Manraj","Masi","Carl")
$task  = @("Demo","Inperson Technical Assistance","Lockup","Operator","Other","Pickup Large PA","Pickup Mic","Pickup PC","Pickup Projector","Pickup Skype Kit","Pickup Small PA","Proactive Classroom Check","Replace Battery","SCLD Student Event","SCLD Student Logout","Setup Large PA","Setup Mic","Setup PC","Setup Projector","Setup Skype Kit","Setup Small PA",$null)
$bldg  = @("BC","BCS","BRG","BSB","CB","CC","CFA","CFT","CLH","CSQ","DB","ELC","FC","FRQ","HNE","K","KT","LAS","LSB","LUM","MC","ATK")

for ($i = 0; $i -lt $staff.Length; $i++) {
    $r = 4 + $i
    $db.Cells.Item($r, 1).Value = $staff[$i]
    $db.Cells.Item($r, 2).Value = $task[$i]
    $db.Cells.Item($r, 3).Value = $bldg[$i]
}

# revert the now-unused trailing rows back to the sheet's default row height
$db.Rows("40").AutoFit()
$db.Rows("41").AutoFit()
$db.Rows("42").AutoFit()
$db.Rows("47").AutoFit()
$db.Rows("48").AutoFit()
$db.Rows("49").AutoFit()
$db.Rows("50").AutoFit()
$db.Rows("51").AutoFit()
$db.Rows("52").AutoFit()

# ---------------------------------------------------------------------------
# 4. Selections - scroll/select to match where Jeannine left off, and make
#    sure "Logs" ends up as the active sheet again.
# ---------------------------------------------------------------------------
$db.Range("B2:B25").Select()

$logs.Activate()
$excel.ActiveWindow.ScrollRow = 37
$logs.Range("F61").Select()
